$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new inventory rows (21-23)
$ws.Range("A21").Value2 = "Pcb with components"
$ws.Range("B21").Value2 = 1
$ws.Range("C21").Value2 = 68
$ws.Range("D21").Value2 = 1

$ws.Range("A22").Value2 = "regulators and stuff"
$ws.Range("B22").Value2 = 1
$ws.Range("C22").Value2 = 15.6
$ws.Range("D22").Value2 = 1

$ws.Range("A23").Value2 = "acid"
$ws.Range("B23").Value2 = 1
$ws.Range("C23").Value2 = 20

# Update the "Used" array formula to include the new parts' cost
$ws.Range("G11").FormulaArray = "=SUM(C2:C12*D2:D12) +E20+C21+C22+C23"

# Recalculate so dependent cells (G11, G12) pick up new values
$wb.Application.Calculate()

# Update the selected cell to reflect where the author ended up
$ws.Range("G12").Select()
